{"js": "// Replace each two-digit multiplication equation in the document with its\n// updated counterpart, per the commit's regenerated answer set.\nconst replacements = [\n  [\"22\u00d727=594\", \"95\u00d779=7505\"],\n  [\"59\u00d728=1652\", \"26\u00d739=1014\"],\n  [\"62\u00d725=1550\", \"58\u00d764=3712\"],\n  [\"60\u00d728=1680\", \"28\u00d784=2352\"],\n  [\"64\u00d755=3520\", \"66\u00d752=3432\"],\n  [\"62\u00d731=1922\", \"80\u00d743=3440\"],\n  [\"81\u00d734=2754\", \"36\u00d763=2268\"],\n  [\"82\u00d783=6806\", \"27\u00d744=1188\"],\n  [\"52\u00d754=2808\", \"80\u00d720=1600\"],\n  [\"45\u00d724=1080\", \"95\u00d764=6080\"],\n  [\"46\u00d739=1794\", \"30\u00d751=1530\"],\n  [\"34\u00d751=1734\", \"82\u00d767=5494\"],\n  [\"63\u00d762=3906\", \"28\u00d767=1876\"],\n  [\"60\u00d780=4800\", \"57\u00d754=3078\"],\n  [\"48\u00d741=1968\", \"37\u00d737=1369\"],\n  [\"29\u00d765=1885\", \"18\u00d784=1512\"],\n  [\"69\u00d749=3381\", \"56\u00d750=2800\"],\n  [\"95\u00d721=1995\", \"27\u00d733=891\"],\n  [\"36\u00d760=2160\", \"31\u00d733=1023\"],\n  [\"46\u00d741=1886\", \"50\u00d780=4000\"],\n  [\"80\u00d788=7040\", \"29\u00d762=1798\"],\n  [\"32\u00d793=2976\", \"25\u00d753=1325\"],\n  [\"31\u00d764=1984\", \"44\u00d732=1408\"],\n  [\"60\u00d716=960\", \"75\u00d760=4500\"],\n  [\"42\u00d799=4158\", \"41\u00d746=1886\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation with its regenerated\n# counterpart, per the commit's updated answer set.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"22\u00d727=594\", \"95\u00d779=7505\"),\n    @(\"59\u00d728=1652\", \"26\u00d739=1014\"),\n    @(\"62\u00d725=1550\", \"58\u00d764=3712\"),\n    @(\"60\u00d728=1680\", \"28\u00d784=2352\"),\n    @(\"64\u00d755=3520\", \"66\u00d752=3432\"),\n    @(\"62\u00d731=1922\", \"80\u00d743=3440\"),\n    @(\"81\u00d734=2754\", \"36\u00d763=2268\"),\n    @(\"82\u00d783=6806\", \"27\u00d744=1188\"),\n    @(\"52\u00d754=2808\", \"80\u00d720=1600\"),\n    @(\"45\u00d724=1080\", \"95\u00d764=6080\"),\n    @(\"46\u00d739=1794\", \"30\u00d751=1530\"),\n    @(\"34\u00d751=1734\", \"82\u00d767=5494\"),\n    @(\"63\u00d762=3906\", \"28\u00d767=1876\"),\n    @(\"60\u00d780=4800\", \"57\u00d754=3078\"),\n    @(\"48\u00d741=1968\", \"37\u00d737=1369\"),\n    @(\"29\u00d765=1885\", \"18\u00d784=1512\"),\n    @(\"69\u00d749=3381\", \"56\u00d750=2800\"),\n    @(\"95\u00d721=1995\", \"27\u00d733=891\"),\n    @(\"36\u00d760=2160\", \"31\u00d733=1023\"),\n    @(\"46\u00d741=1886\", \"50\u00d780=4000\"),\n    @(\"80\u00d788=7040\", \"29\u00d762=1798\"),\n    @(\"32\u00d793=2976\", \"25\u00d753=1325\"),\n    @(\"31\u00d764=1984\", \"44\u00d732=1408\"),\n    @(\"60\u00d716=960\", \"75\u00d760=4500\"),\n    @(\"42\u00d799=4158\", \"41\u00d746=1886\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2) | Out-Null\n}\n"}
